# "1st changes of mifos to finflux"
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet - this pushes the old N/O/P ("Late"/"Outstanding"/"Original")
# columns one slot to the right (O/P/Q), and makes the "Repayment
# schedule" tab the active sheet/tab instead of "Transactions".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column in front of column N; everything from N onward
# (values + styles) shifts right automatically.
$ws.Columns("N").Insert()

# The newly inserted column picks up the width of its left neighbour
# (column M), matching Excel's normal "insert column" behaviour.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab (was "Transactions").
$ws.Activate()
$ws.Range("S7").Select()
